$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(122, 2).Value = 6720873
$ws.Cells.Item(122, 6).Value = 'Sportivo Luqueno'
$ws.Cells.Item(122, 7).Value = 'Sportivo Trinidense'
$ws.Cells.Item(122, 9).Value = 2
$ws.Cells.Item(122, 10).Value = 'D'
$ws.Cells.Item(122, 11).Value = 2.625
$ws.Cells.Item(122, 12).Value = 3.1
$ws.Cells.Item(122, 13).Value = 2.5
$ws.Cells.Item(122, 14).Value = 2.3
$ws.Cells.Item(122, 15).Value = 3.1
$ws.Cells.Item(122, 16).Value = 2.9
$ws.Cells.Item(122, 17).Value = -0.25
$ws.Cells.Item(122, 18).Value = 2.025
$ws.Cells.Item(122, 19).Value = 1.775
$ws.Cells.Item(122, 21).Value = 1.95
$ws.Cells.Item(122, 22).Value = 1.85
$ws.Cells.Item(122, 23).Value = -1
$ws.Cells.Item(122, 24).Value = 2.1
$ws.Cells.Item(122, 26).Value = -0.5
$ws.Cells.Item(122, 27).Value = 0.3875
$ws.Cells.Item(122, 28).Value = 0.95
$ws.Cells.Item(122, 29).Value = -1
$ws.Cells.Item(124, 2).Value = 6720843
$ws.Cells.Item(124, 6).Value = 'Cerro Porteno'
$ws.Cells.Item(124, 7).Value = 'Libertad Asuncion'
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 'H'
$ws.Cells.Item(124, 11).Value = 2.375
$ws.Cells.Item(124, 12).Value = 3.2
$ws.Cells.Item(124, 13).Value = 2.7
$ws.Cells.Item(124, 14).Value = 3.75
$ws.Cells.Item(124, 15).Value = 3.3
$ws.Cells.Item(124, 16).Value = 1.85
$ws.Cells.Item(124, 17).Value = 0.5
$ws.Cells.Item(124, 18).Value = 1.9
$ws.Cells.Item(124, 19).Value = 1.9
$ws.Cells.Item(124, 21).Value = 1.925
$ws.Cells.Item(124, 22).Value = 1.875
$ws.Cells.Item(124, 23).Value = 2.75
$ws.Cells.Item(124, 24).Value = -1
$ws.Cells.Item(124, 26).Value = 0.8999999999999999
$ws.Cells.Item(124, 27).Value = -1
$ws.Cells.Item(124, 28).Value = -1
$ws.Cells.Item(124, 29).Value = 0.875
$ws.Cells.Item(260, 2).Value = 7493431
$ws.Cells.Item(260, 6).Value = 'Sportivo Trinidense'
$ws.Cells.Item(260, 7).Value = 'Guairena FC'
$ws.Cells.Item(260, 8).Value = 7
$ws.Cells.Item(260, 10).Value = 'H'
$ws.Cells.Item(260, 11).Value = 2.05
$ws.Cells.Item(260, 12).Value = 3.3
$ws.Cells.Item(260, 13).Value = 3.3
$ws.Cells.Item(260, 14).Value = 2.6
$ws.Cells.Item(260, 15).Value = 3.1
$ws.Cells.Item(260, 16).Value = 2.6
$ws.Cells.Item(260, 17).Value = 0
$ws.Cells.Item(260, 18).Value = 1.925
$ws.Cells.Item(260, 19).Value = 1.875
$ws.Cells.Item(260, 20).Value = 2.5
$ws.Cells.Item(260, 21).Value = 2
$ws.Cells.Item(260, 22).Value = 1.8
$ws.Cells.Item(260, 23).Value = 1.6
$ws.Cells.Item(260, 25).Value = -1
$ws.Cells.Item(260, 26).Value = 0.925
$ws.Cells.Item(260, 27).Value = -1
$ws.Cells.Item(260, 28).Value = 1
$ws.Cells.Item(260, 29).Value = -1
$ws.Cells.Item(261, 2).Value = 7493310
$ws.Cells.Item(261, 6).Value = 'Libertad Asuncion'
$ws.Cells.Item(261, 7).Value = 'Tacuary'
$ws.Cells.Item(261, 8).Value = 1
$ws.Cells.Item(261, 10).Value = 'A'
$ws.Cells.Item(261, 11).Value = 1.363
$ws.Cells.Item(261, 12).Value = 5
$ws.Cells.Item(261, 13).Value = 7
$ws.Cells.Item(261, 14).Value = 1.571
$ws.Cells.Item(261, 15).Value = 4.2
$ws.Cells.Item(261, 16).Value = 4.75
$ws.Cells.Item(261, 17).Value = -0.75
$ws.Cells.Item(261, 18).Value = 1.8
$ws.Cells.Item(261, 19).Value = 2
$ws.Cells.Item(261, 20).Value = 2.75
$ws.Cells.Item(261, 21).Value = 1.8
$ws.Cells.Item(261, 22).Value = 2
$ws.Cells.Item(261, 23).Value = -1
$ws.Cells.Item(261, 25).Value = 3.75
$ws.Cells.Item(261, 26).Value = -1
$ws.Cells.Item(261, 27).Value = 1
$ws.Cells.Item(261, 28).Value = 0.4
$ws.Cells.Item(261, 29).Value = -0.5
$ws.Cells.Item(263, 2).Value = 7493312
$ws.Cells.Item(263, 6).Value = 'Cerro Porteno'
$ws.Cells.Item(263, 7).Value = 'Guarani Asuncion'
$ws.Cells.Item(263, 8).Value = 4
$ws.Cells.Item(263, 9).Value = 0
$ws.Cells.Item(263, 10).Value = 'H'
$ws.Cells.Item(263, 11).Value = 1.7
$ws.Cells.Item(263, 12).Value = 3.6
$ws.Cells.Item(263, 13).Value = 4.333
$ws.Cells.Item(263, 14).Value = 1.727
$ws.Cells.Item(263, 15).Value = 3.75
$ws.Cells.Item(263, 16).Value = 4.2
$ws.Cells.Item(263, 17).Value = -0.5
$ws.Cells.Item(263, 18).Value = 1.8
$ws.Cells.Item(263, 19).Value = 2
$ws.Cells.Item(263, 20).Value = 2.75
$ws.Cells.Item(263, 21).Value = 1.875
$ws.Cells.Item(263, 22).Value = 1.925
$ws.Cells.Item(263, 23).Value = 0.7270000000000001
$ws.Cells.Item(263, 24).Value = -1
$ws.Cells.Item(263, 26).Value = 0.8
$ws.Cells.Item(263, 27).Value = -1
$ws.Cells.Item(263, 28).Value = 0.875
$ws.Cells.Item(263, 29).Value = -1
$ws.Cells.Item(264, 2).Value = 7493433
$ws.Cells.Item(264, 6).Value = 'Sportivo Luqueno'
$ws.Cells.Item(264, 7).Value = 'Nacional Asuncion'
$ws.Cells.Item(264, 8).Value = 1
$ws.Cells.Item(264, 10).Value = 'D'
$ws.Cells.Item(264, 11).Value = 2.75
$ws.Cells.Item(264, 12).Value = 3.2
$ws.Cells.Item(264, 13).Value = 2.4
$ws.Cells.Item(264, 14).Value = 2.75
$ws.Cells.Item(264, 15).Value = 3.1
$ws.Cells.Item(264, 16).Value = 2.45
$ws.Cells.Item(264, 18).Value = 1.75
$ws.Cells.Item(264, 19).Value = 2.05
$ws.Cells.Item(264, 21).Value = 2
$ws.Cells.Item(264, 22).Value = 1.8
$ws.Cells.Item(264, 24).Value = 2.1
$ws.Cells.Item(264, 25).Value = -1
$ws.Cells.Item(264, 26).Value = 0.375
$ws.Cells.Item(264, 27).Value = -0.5
$ws.Cells.Item(264, 28).Value = -0.5
$ws.Cells.Item(264, 29).Value = 0.4
$ws.Cells.Item(265, 2).Value = 7493311
$ws.Cells.Item(265, 6).Value = 'General Caballero JLM'
$ws.Cells.Item(265, 7).Value = 'Olimpia Asuncion'
$ws.Cells.Item(265, 8).Value = 0
$ws.Cells.Item(265, 9).Value = 1
$ws.Cells.Item(265, 10).Value = 'A'
$ws.Cells.Item(265, 11).Value = 3.4
$ws.Cells.Item(265, 12).Value = 3.3
$ws.Cells.Item(265, 13).Value = 2
$ws.Cells.Item(265, 14).Value = 3.2
$ws.Cells.Item(265, 15).Value = 3.25
$ws.Cells.Item(265, 16).Value = 2.1
$ws.Cells.Item(265, 17).Value = 0.25
$ws.Cells.Item(265, 18).Value = 1.95
$ws.Cells.Item(265, 19).Value = 1.85
$ws.Cells.Item(265, 20).Value = 2.25
$ws.Cells.Item(265, 21).Value = 1.775
$ws.Cells.Item(265, 22).Value = 2.025
$ws.Cells.Item(265, 23).Value = -1
$ws.Cells.Item(265, 25).Value = 1.1
$ws.Cells.Item(265, 26).Value = -1
$ws.Cells.Item(265, 27).Value = 0.8500000000000001
$ws.Cells.Item(265, 28).Value = -1
$ws.Cells.Item(265, 29).Value = 1.025
$ws.Cells.Item(302, 8).Value = 2
$ws.Cells.Item(302, 9).Value = 0
$ws.Cells.Item(302, 10).Value = 'H'
$ws.Cells.Item(302, 14).Value = 2.25
$ws.Cells.Item(302, 15).Value = 3.25
$ws.Cells.Item(302, 18).Value = 2.025
$ws.Cells.Item(302, 19).Value = 1.775
$ws.Cells.Item(302, 21).Value = 2.025
$ws.Cells.Item(302, 22).Value = 1.775
$ws.Cells.Item(302, 23).Value = 1.25
$ws.Cells.Item(302, 24).Value = -1
$ws.Cells.Item(302, 25).Value = -1
$ws.Cells.Item(302, 26).Value = 1.025
$ws.Cells.Item(302, 27).Value = -1
$ws.Cells.Item(302, 28).Value = -0.5
$ws.Cells.Item(302, 29).Value = 0.3875
$ws.Cells.Item(303, 2).Value = 7609135
$ws.Cells.Item(303, 5).Value = 45347.76041666666
$ws.Cells.Item(303, 6).Value = 'Sportivo Trinidense'
$ws.Cells.Item(303, 7).Value = 'Sportivo Ameliano'
$ws.Cells.Item(303, 11).Value = 2.7
$ws.Cells.Item(303, 12).Value = 3.2
$ws.Cells.Item(303, 13).Value = 2.4
$ws.Cells.Item(303, 14).Value = 2.625
$ws.Cells.Item(303, 15).Value = 3.2
$ws.Cells.Item(303, 16).Value = 2.45
$ws.Cells.Item(303, 17).Value = 0
$ws.Cells.Item(303, 18).Value = 2
$ws.Cells.Item(303, 19).Value = 1.8
$ws.Cells.Item(303, 20).Value = 2.5
$ws.Cells.Item(303, 21).Value = 1.925
$ws.Cells.Item(303, 22).Value = 1.875
$ws.Cells.Item(304, 2).Value = 7609137
$ws.Cells.Item(304, 5).Value = 45347.86458333334
$ws.Cells.Item(304, 6).Value = 'Guarani Asuncion'
$ws.Cells.Item(304, 7).Value = 'Olimpia Asuncion'
$ws.Cells.Item(304, 11).Value = 2.6
$ws.Cells.Item(304, 13).Value = 2.5
$ws.Cells.Item(304, 14).Value = 2.8
$ws.Cells.Item(304, 16).Value = 2.3
$ws.Cells.Item(304, 17).Value = 0.25
$ws.Cells.Item(304, 18).Value = 1.75
$ws.Cells.Item(304, 19).Value = 2.05
$ws.Cells.Item(304, 20).Value = 2.25
$ws.Cells.Item(304, 21).Value = 1.825
$ws.Cells.Item(304, 22).Value = 1.975
$ws.Cells.Item(305, 2).Value = 7609192
$ws.Cells.Item(305, 5).Value = 45348.8125
$ws.Cells.Item(305, 6).Value = 'Libertad Asuncion'
$ws.Cells.Item(305, 7).Value = 'Tacuary'
$ws.Cells.Item(305, 11).Value = 1.363
$ws.Cells.Item(305, 12).Value = 4.333
$ws.Cells.Item(305, 13).Value = 7.5
$ws.Cells.Item(305, 14).Value = 1.363
$ws.Cells.Item(305, 15).Value = 4.333
$ws.Cells.Item(305, 16).Value = 7.5
$ws.Cells.Item(305, 17).Value = -1.25
$ws.Cells.Item(305, 18).Value = 1.825
$ws.Cells.Item(305, 19).Value = 1.975
$ws.Cells.Item(305, 20).Value = 2.75
$ws.Cells.Item(305, 21).Value = 1.9
$ws.Cells.Item(305, 22).Value = 1.9

$ws.Rows.Item(306).Delete()

